# RequiredData.xlsx was re-uploaded with refreshed dummy test credentials.
# Update the stored username/password values on the "Credentials" sheet
# (A2 = username, B2 = password) to the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Credentials")

$ws.Range("A2").Value = "mngr608625"
$ws.Range("B2").Value = "uvabujU"
